# Refresh Sargatanas market-board price snapshots + recomputed leve profit figures
# (scheduled runner data update across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 520
$ws.Range("J9").Value = 280
$ws.Range("L9").Value = 280
$ws.Range("N9").Value = -618

$ws.Range("H17").Value = 1292.6229
$ws.Range("J17").Value = 1309.322
$ws.Range("L17").Value = 3927.965999999999
$ws.Range("N17").Value = -4263.965999999999

$ws.Range("H28").Value = 1417.5714
$ws.Range("I28").Value = 1849.25
$ws.Range("K28").Value = 1849.25
$ws.Range("M28").Value = -1364.25

$ws.Range("H129").Value = 1136.5
$ws.Range("I129").Value = 618.4
$ws.Range("K129").Value = 1855.2
$ws.Range("M129").Value = 3144.8

$ws.Range("H137").Value = 2499.1462
$ws.Range("I137").Value = 2184.5356
$ws.Range("K137").Value = 6553.6068
$ws.Range("M137").Value = -4003.6068

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4698.909
$ws.Range("J2").Value = 6615.8335
$ws.Range("L2").Value = 6615.8335
$ws.Range("N2").Value = -6841.8335

$ws.Range("H45").Value = 11290.9
$ws.Range("I45").Value = 3074
$ws.Range("J45").Value = 14812.429
$ws.Range("K45").Value = 3074
$ws.Range("L45").Value = 14812.429
$ws.Range("M45").Value = -2697
$ws.Range("N45").Value = -15566.429

$ws.Range("H61").Value = 5917.921
$ws.Range("I61").Value = 2417.32
$ws.Range("K61").Value = 2417.32
$ws.Range("M61").Value = -2205.32

$ws.Range("H116").Value = 4698.909
$ws.Range("J116").Value = 6615.8335
$ws.Range("L116").Value = 6615.8335
$ws.Range("N116").Value = -11203.8335

$ws.Range("H132").Value = 6074.39
$ws.Range("I132").Value = 4408.3335
$ws.Range("K132").Value = 13225.0005
$ws.Range("M132").Value = -10695.0005

$ws.Range("H136").Value = 5917.921
$ws.Range("I136").Value = 2417.32
$ws.Range("K136").Value = 7251.960000000001
$ws.Range("M136").Value = -4701.960000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4698.909
$ws.Range("J3").Value = 6615.8335
$ws.Range("L3").Value = 6615.8335
$ws.Range("N3").Value = -6843.8335

$ws.Range("H22").Value = 183.27272
$ws.Range("I22").Value = 203.77777
$ws.Range("J22").Value = 91
$ws.Range("K22").Value = 203.77777
$ws.Range("L22").Value = 91
$ws.Range("M22").Value = -30.77777
$ws.Range("N22").Value = -437

$ws.Range("H86").Value = 43525776
$ws.Range("I86").Value = 79613.08
$ws.Range("J86").Value = 100005790
$ws.Range("K86").Value = 79613.08
$ws.Range("L86").Value = 100005790
$ws.Range("M86").Value = -78490.08
$ws.Range("N86").Value = -100008036

$ws.Range("H89").Value = 43525776
$ws.Range("I89").Value = 79613.08
$ws.Range("J89").Value = 100005790
$ws.Range("K89").Value = 398065.4
$ws.Range("L89").Value = 500028950
$ws.Range("M89").Value = -392449.4
$ws.Range("N89").Value = -500040182

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 6107.1665
$ws.Range("I16").Value = 3357.2
$ws.Range("K16").Value = 3357.2
$ws.Range("M16").Value = -3070.2

$ws.Range("H31").Value = 8640.177
$ws.Range("I31").Value = 3791.889
$ws.Range("K31").Value = 3791.889
$ws.Range("M31").Value = -3496.889

$ws.Range("H34").Value = 8640.177
$ws.Range("I34").Value = 3791.889
$ws.Range("K34").Value = 3791.889
$ws.Range("M34").Value = -3589.889

$ws.Range("H86").Value = 4469696.5
$ws.Range("I86").Value = 8933779
$ws.Range("J86").Value = 5614.857
$ws.Range("K86").Value = 8933779
$ws.Range("L86").Value = 5614.857
$ws.Range("M86").Value = -8932656
$ws.Range("N86").Value = -7860.857

$ws.Range("H89").Value = 4469696.5
$ws.Range("I89").Value = 8933779
$ws.Range("J89").Value = 5614.857
$ws.Range("K89").Value = 44668895
$ws.Range("L89").Value = 28074.285
$ws.Range("M89").Value = -44663279
$ws.Range("N89").Value = -39306.285

$ws.Range("H99").Value = 7673.5
$ws.Range("I99").Value = 8773.777
$ws.Range("J99").Value = 6773.273
$ws.Range("K99").Value = 8773.777
$ws.Range("L99").Value = 6773.273
$ws.Range("M99").Value = -7275.777
$ws.Range("N99").Value = -9769.273000000001

$ws.Range("H107").Value = 2068.64
$ws.Range("I107").Value = 1525.9333
$ws.Range("J107").Value = 2882.7
$ws.Range("K107").Value = 1525.9333
$ws.Range("L107").Value = 2882.7
$ws.Range("M107").Value = 394.0667000000001
$ws.Range("N107").Value = -6722.7

$ws.Range("H113").Value = 6107.1665
$ws.Range("I113").Value = 3357.2
$ws.Range("K113").Value = 3357.2
$ws.Range("M113").Value = -1187.2

$ws.Range("H126").Value = 7673.5
$ws.Range("I126").Value = 8773.777
$ws.Range("J126").Value = 6773.273
$ws.Range("K126").Value = 26321.331
$ws.Range("L126").Value = 20319.819
$ws.Range("M126").Value = -23851.331
$ws.Range("N126").Value = -25259.819

$ws.Range("H132").Value = 6232.1763
$ws.Range("I132").Value = 1810.9333
$ws.Range("J132").Value = 9722.632
$ws.Range("K132").Value = 5432.7999
$ws.Range("L132").Value = 29167.896
$ws.Range("M132").Value = -2902.7999
$ws.Range("N132").Value = -34227.896

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3334971.5
$ws.Range("J5").Value = 4125
$ws.Range("L5").Value = 12375
$ws.Range("N5").Value = -12599

$ws.Range("H97").Value = 1099.5
$ws.Range("I97").Value = 799
$ws.Range("J97").Value = 1400
$ws.Range("K97").Value = 2397
$ws.Range("L97").Value = 4200
$ws.Range("M97").Value = -1901
$ws.Range("N97").Value = -5192

$ws.Range("H135").Value = 3334971.5
$ws.Range("J135").Value = 4125
$ws.Range("L135").Value = 37125
$ws.Range("N135").Value = -42195

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8181.6
$ws.Range("I70").Value = 7115.4287
$ws.Range("J70").Value = 10669.333
$ws.Range("K70").Value = 7115.4287
$ws.Range("L70").Value = 10669.333
$ws.Range("M70").Value = -6845.4287
$ws.Range("N70").Value = -11209.333

$ws.Range("H73").Value = 8181.6
$ws.Range("I73").Value = 7115.4287
$ws.Range("J73").Value = 10669.333
$ws.Range("K73").Value = 7115.4287
$ws.Range("L73").Value = 10669.333
$ws.Range("M73").Value = -6179.4287
$ws.Range("N73").Value = -12541.333

$ws.Range("H97").Value = 975.2727
$ws.Range("I97").Value = 841
$ws.Range("J97").Value = 1333.3334
$ws.Range("K97").Value = 841
$ws.Range("L97").Value = 1333.3334
$ws.Range("M97").Value = -345
$ws.Range("N97").Value = -2325.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0

$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0

$ws.Range("H122").Value = 7666.8887
$ws.Range("I122").Value = 4999
$ws.Range("K122").Value = 14997
$ws.Range("M122").Value = -12547

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 5552.7144
$ws.Range("I126").Value = 2990
$ws.Range("J126").Value = 7474.75
$ws.Range("K126").Value = 8970
$ws.Range("L126").Value = 22424.25
$ws.Range("M126").Value = -6500
$ws.Range("N126").Value = -27364.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("N69").ClearContents()
$ws.Range("N72").ClearContents()
